# Update the public EPEX Spot prices workbook:
#  - "Prix Spot": add a new day column CN ("13-sep") with its 24 hourly values.
#  - "Gaz" and "CO2": append a new date row (2025-09-11) with its last price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add column CN (92) mirroring the style of column CM (91)
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold, centered, bordered) from CM1 to
# CN1 so the new header cell keeps the same style as the rest of the header
# row, then set its text.
$wsPrix.Cells.Item(1, 91).Copy()
$wsPrix.Cells.Item(1, 92).PasteSpecial(-4122)   # xlPasteFormats
$wsPrix.Cells.Item(1, 92).Value = "13-sep"

# New hourly values for 13-sep (rows 2..25 correspond to the 24 hourly slots)
$wsPrix.Cells.Item(2, 92).Value = 9.45
$wsPrix.Cells.Item(3, 92).Value = 9.17
$wsPrix.Cells.Item(4, 92).Value = 10.35
$wsPrix.Cells.Item(5, 92).Value = 7.39
$wsPrix.Cells.Item(6, 92).Value = 7
$wsPrix.Cells.Item(7, 92).Value = 7.75
$wsPrix.Cells.Item(8, 92).Value = 10.55
$wsPrix.Cells.Item(9, 92).Value = 9.63
$wsPrix.Cells.Item(10, 92).Value = 13.63
$wsPrix.Cells.Item(11, 92).Value = 12.42
$wsPrix.Cells.Item(12, 92).Value = 9.62
$wsPrix.Cells.Item(13, 92).Value = 2.37
$wsPrix.Cells.Item(14, 92).Value = 0
$wsPrix.Cells.Item(15, 92).Value = -0.01
$wsPrix.Cells.Item(16, 92).Value = -0.01
$wsPrix.Cells.Item(17, 92).Value = 0
$wsPrix.Cells.Item(18, 92).Value = 0
$wsPrix.Cells.Item(19, 92).Value = 0.65
$wsPrix.Cells.Item(20, 92).Value = 14
$wsPrix.Cells.Item(21, 92).Value = 17.68
$wsPrix.Cells.Item(22, 92).Value = 25.27
$wsPrix.Cells.Item(23, 92).Value = 10.51
$wsPrix.Cells.Item(24, 92).Value = 29.45
$wsPrix.Cells.Item(25, 92).Value = 25.73

# ---------------------------------------------------------------------------
# Helper: write a date-like string (e.g. "2025-09-11") into a cell as plain
# text rather than letting it be auto-recognised and converted into a date
# serial number. We build the string via a throw-away formula cell (which
# keeps it typed as text) and copy only the resulting value across, so the
# destination cell's number format/style stays completely untouched.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $text) {
    $scratch = $ws.Cells.Item(1048576, 16384)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $target = $ws.Cells.Item($row, $col)
    $target.PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 89 (2025-09-11, 31.8)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
Set-TextValue $wsGaz 89 1 "2025-09-11"
$wsGaz.Cells.Item(89, 2).Value = 31.8

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 89 (2025-09-11, 75.25)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
Set-TextValue $wsCo2 89 1 "2025-09-11"
$wsCo2.Cells.Item(89, 2).Value = 75.25
